# "Adding Custom Annotations in Framework"
#
# RUNMANAGER sheet (sheet1): drop the Priority/Count columns (D, E), flip the
# "execute" flag for the two existing tests to "no", and append a new
# "leaveFeatureTest" row that is flagged to execute.
#
# DATA sheet (sheet2): append the matching test-data row for the new
# "leaveFeatureTest" test case.
#
# Also update the active sheet / selections: RUNMANAGER becomes the active
# tab (selection on C2) and DATA's selection grows to A5:E12.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RUNMANAGER
$ws2 = $wb.Worksheets.Item(2)   # DATA

# --- RUNMANAGER: remove the trailing Priority (D) / Count (E) columns ---
$ws1.Columns.Item(4).EntireColumn.Delete() | Out-Null
$ws1.Columns.Item(4).EntireColumn.Delete() | Out-Null

# Existing rows now execute = "no"
$ws1.Range("C2").Value = "no"
$ws1.Range("C3").Value = "no"

# New row for the leave-feature test, flagged to execute
$ws1.Range("A4").Value = "leaveFeatureTest"
$ws1.Range("B4").Value = "To test Leave Feature"
$ws1.Range("C4").Value = "yes"

# --- DATA: add the test-data row for the new test case ---
$ws2.Range("A4").Value = "leaveFeatureTest"
$ws2.Range("B4").Value = "yes"
$ws2.Range("C4").Value = "chrome"
$ws2.Range("D4").Value = "Admin"
$ws2.Range("E4").Value = "admin123"

# --- Selections / active sheet ---
$ws2.Activate()
$ws2.Range("A5:E12").Select() | Out-Null

$ws1.Activate()
$ws1.Range("C2").Select() | Out-Null
